$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("B26").Value = "OPQA-512"
$ws.Range("B27").Value = "OPQA-516"
$ws.Range("B28").Value = "OPQA-517"
$ws.Range("B29").Value = "OPQA-518"

$ws.Range("B32").Select()
